# Swap the contents of columns C ("TYPE") and D ("STATE") across every
# used row of the active sheet (header row included). The sheet shipped
# with these two columns' data transposed relative to their headers;
# this restores TYPE values to column C and STATE values to column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($i = 1; $i -le $lastRow; $i++) {
    $cCell = $ws.Cells.Item($i, 3)
    $dCell = $ws.Cells.Item($i, 4)

    $cVal = $cCell.Value()
    $dVal = $dCell.Value()

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
